# Swap the data rows for "Net Profit" (row 9) and "Pretax Margin" (row 11).
# Row 10 ("Operating Profit Margin") stays where it is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (before-edit) values of row 9 and row 11
$row9 = @($ws.Range("A9").Value2, $ws.Range("B9").Value2, $ws.Range("C9").Value2, $ws.Range("D9").Value2)
$row11 = @($ws.Range("A11").Value2, $ws.Range("B11").Value2, $ws.Range("C11").Value2, $ws.Range("D11").Value2)

# Write row 11's original content into row 9
$ws.Range("A9").Value = $row11[0]
$ws.Range("B9").Value = $row11[1]
$ws.Range("C9").Value = $row11[2]
$ws.Range("D9").Value = $row11[3]

# Write row 9's original content into row 11
$ws.Range("A11").Value = $row9[0]
$ws.Range("B11").Value = $row9[1]
$ws.Range("C11").Value = $row9[2]
$ws.Range("D11").Value = $row9[3]
